$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row of data (row 13) continuing the existing table
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 2
